$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions) - simple "want to go" count bumps
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 2630
$ws1.Range("F8").Value  = 54
$ws1.Range("F11").Value = 1530
$ws1.Range("F12").Value = 4
$ws1.Range("F14").Value = 630
$ws1.Range("F16").Value = 1359
$ws1.Range("F17").Value = 18
$ws1.Range("F18").Value = 538
$ws1.Range("F19").Value = 3807
$ws1.Range("F21").Value = 3299
$ws1.Range("F23").Value = 21
$ws1.Range("F24").Value = 2188
$ws1.Range("F26").Value = 304
$ws1.Range("F28").Value = 23
$ws1.Range("F29").Value = 1162
$ws1.Range("F30").Value = 762
$ws1.Range("F32").Value = 1038
$ws1.Range("F33").Value = 1032

# ---------------------------------------------------------------
# Sheet 2: 演出 (Shows) - simple "want to go" count bumps
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 101
$ws2.Range("F18").Value = 248
$ws2.Range("F19").Value = 193
$ws2.Range("F20").Value = 476

# ---------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life) - value bump plus row restructure:
#   * F2 bumps 244 -> 243
#   * Old row 4 ("Paradox Live" event, now expired) is removed
#   * Old row 5 (EVANGELION event) becomes row 4, with its "want to
#     go" count updated from 528 to 529 and ticket price now a
#     numeric 20 (was a non-numeric "不可售" placeholder)
#   * Old row 6 (Crayon Shin-chan event) becomes row 5 unchanged
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 243

# B4/B5 hold plain "yyyy-mm-dd" text (matching the rest of the column),
# not real dates - force text entry so Excel doesn't coerce them into a
# date serial, then drop the number format back off the cell.
$ws3.Range("B4").NumberFormat = "@"
$ws3.Range("B4").Value = "2024-10-15"
$ws3.Range("B4").Style = "Normal"
$ws3.Range("C4").Value = "北京·EVANGELION× PrismLand · 新世纪福音战士官方授权主题店"
$ws3.Range("E4").Value = "2024.10.15 00:00-12.15 23:59"
$ws3.Range("F4").Value = 529
$ws3.Range("G4").Value = 20
$ws3.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=93039"
$ws3.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202409/n32CfRya1727584778969.jpeg"

$ws3.Range("B5").NumberFormat = "@"
$ws3.Range("B5").Value = "2024-10-17"
$ws3.Range("B5").Style = "Normal"
$ws3.Range("C5").Value = "北京·蜡笔小新：我们的恐龙日记x HAPPY ZOO 主题咖啡厅"
$ws3.Range("E5").Value = "2024.10.17 00:00-10.27 23:59"
$ws3.Range("F5").Value = 103
$ws3.Range("G5").Value = 10
$ws3.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=93224"
$ws3.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202410/nzGP5KRA1728526131597.png"

$ws3.Rows.Item(6).Delete()

# ---------------------------------------------------------------
# Sheet 4: 全部类型 (All Types) - combined view, value bumps only
# (this sheet never included the removed "Paradox Live" row, so no
# structural change is needed here)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 243
$ws4.Range("F9").Value  = 529
$ws4.Range("F11").Value = 2630
$ws4.Range("F12").Value = 2630
$ws4.Range("F22").Value = 1530
$ws4.Range("F26").Value = 101
$ws4.Range("F27").Value = 1359
$ws4.Range("F28").Value = 18
$ws4.Range("F29").Value = 538
$ws4.Range("F31").Value = 3807
$ws4.Range("F33").Value = 3299
$ws4.Range("F35").Value = 2188
$ws4.Range("F37").Value = 304
$ws4.Range("F39").Value = 23
$ws4.Range("F40").Value = 1162
$ws4.Range("F43").Value = 248
$ws4.Range("F44").Value = 193
$ws4.Range("F45").Value = 476
$ws4.Range("F46").Value = 762
$ws4.Range("F48").Value = 1038
$ws4.Range("F49").Value = 1032

$wb.Save()
